# Updated symbol list on Wed Dec 21 21:57:13 UTC 2022 with GitHub Actions
#
# This re-applies the latest price/volume scrape for the cryptos worksheet:
# most rows keep their coin identity but get a refreshed Price (column D);
# rows 18-24 additionally had their coin ranking shift by one position
# (each coin's row moved down one slot, with "One" rotating back up into
# row 18 with a freshly scraped price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Column D stores numeric-looking prices (e.g. "246.67") as plain text,
    # not numbers. A leading apostrophe forces Excel to keep the value as
    # text instead of silently converting it to a real number; resetting
    # the cell style to "Normal" afterwards drops the transient
    # quote-prefix style Excel applies so the cell stays style-less, just
    # like the other text cells on this sheet.
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# --- Straightforward price refreshes (coin identity unchanged) ---
Set-TextCell "D2" "246.98"
Set-TextCell "D3" "22.48"
Set-TextCell "D4" "5.248"
Set-TextCell "D5" "0.05690"
Set-TextCell "D6" "3.417"
Set-TextCell "D7" "6.297"
Set-TextCell "D9" "0.8602"
Set-TextCell "D10" "0.1411"
Set-TextCell "D11" "0.07337"
Set-TextCell "D12" "0.03048"
Set-TextCell "D14" "0.09388"
Set-TextCell "D15" "3.867"
Set-TextCell "D16" "0.001582"
Set-TextCell "D17" "0.04788"

# --- Rows 18-24: coin ranking rotated by one slot ---
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D18" "0.0005842"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D19" "0.006417"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell "D20" "0.005026"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell "D21" "0.0009969"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell "D22" "0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D23" "3.693"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D24" "2.193"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- More straightforward price refreshes further down the sheet ---
Set-TextCell "D26" "0.1332"
Set-TextCell "D40" "0.03915"
Set-TextCell "D41" "0.006775"
Set-TextCell "D42" "0.1067"
Set-TextCell "D43" "0.002670"
Set-TextCell "D44" "0.007539"
Set-TextCell "D45" "0.00005596"
Set-TextCell "D48" "0.1966"
